# RinominaPdf/rinomina.xlsx: strip the ".pdf" suffix from the sample rename
# values and update the header label ("nuovo" -> "new").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1 header: "old" stays the same, "nuovo" becomes "new"
$ws.Range("B1").Value = "new"

# Row 2 sample values: drop the ".pdf" extension
$ws.Range("A2").Value = "old1"

# Row 3 sample values: drop the ".pdf" extension
# (B3 is written before B2 to reproduce the original authoring order)
$ws.Range("B3").Value = "nuovo2"
$ws.Range("B2").Value = "nuovo1"
$ws.Range("A3").Value = "old2"

# Move/keep the active selection on A3
$ws.Range("A3").Select()
